$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Golang Developer", "https://www.dice.com/job-detail/23964ab8-0c28-462c-92c5-977132d21423", "Remote or Hybrid in Cleveland, Ohio", "Contract", "DOE", "Crea Services LLC"),
    @("Golang Architect / Principal Backend Architect", "https://www.dice.com/job-detail/eb2d20e5-cb44-468e-ac55-a2185bcb3be9", "Atlanta, Georgia", "Contract, Third Party", "`$DOE", "Stellite Works LLC"),
    @("Golang Engineer", "https://www.dice.com/job-detail/aef73333-4197-4d3d-b077-2e73a54f8988", "Remote", "Contract", "50 - 55", "Arnex Solutions LLC"),
    @("Software Engineer Backend", "https://www.dice.com/job-detail/0519334e-e679-4fd9-b986-1306a3644194", "Fremont, California", "Contract, Third Party", "Depends on Experience", "Tekfortune Inc.")
)

$startRow = 35
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$c]
    }
}
